# Update DEQM Capability Statement Producer Server workbook to R4 / 1.1.0
$wb = $excel.ActiveWorkbook

# --- "meta" sheet: bump version / fhirVersion / ig url ---
$wsMeta = $wb.Worksheets.Item("meta")
$wsMeta.Range("B3").Value = "1.1.0"
$wsMeta.Range("B4").Value = "4.0.0"
$wsMeta.Range("B6").Value = "http://hl7.org/fhir/us/davinci-deqm/ImplementationGuide/hl7.fhir.us.davinci-deqm-1.1.0"

# --- "profiles" sheet: drop the "/STU3" path segment from each Profile url ---
$wsProfiles = $wb.Worksheets.Item("profiles")
$wsProfiles.Range("A2").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/devicerequest-deqm"
$wsProfiles.Range("A3").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/medicationadministration-deqm"
$wsProfiles.Range("A4").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/summary-measurereport-deqm"
$wsProfiles.Range("A5").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/datax-measurereport-deqm"
$wsProfiles.Range("A6").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/practitioner-deqm"
$wsProfiles.Range("A7").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/deviceusestatement-deqm"
$wsProfiles.Range("A8").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/organization-deqm"
$wsProfiles.Range("A9").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/coverage-deqm"
$wsProfiles.Range("A10").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/indv-measurereport-deqm"
$wsProfiles.Range("A11").Value = "http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/medicationrequest-deqm"
$wsProfiles.Range("A12").Value = "!http://hl7.org/fhir/us/davinci-deqm/StructureDefinition/subscription-deqm"

# --- selections / active sheet ---
# meta sheet cursor moves to D5 (not the active tab)
$wsMeta.Activate()
$wsMeta.Range("D5").Select()

# profiles sheet becomes the active tab with cursor at A14
$wsProfiles.Activate()
$wsProfiles.Range("A14").Select()
